$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry holds the new D (Price) / E (Volume 1h) text for a row that
# changed in this data refresh. A $null value means that column did not
# change for that row.
$updates = @(
    @{Row=2; D="24.793.98"; E="  +0.66%  "},
    @{Row=3; D="1.702.39"; E="  +0.04%  "},
    @{Row=4; D="1.003"; E=$null},
    @{Row=5; D="316.97"; E="  +0.60%  "},
    @{Row=6; D=$null; E="  +0.29%  "},
    @{Row=7; D="0.3931"; E="  -0.49%  "},
    @{Row=8; D="0.4043"; E="  -0.03%  "},
    @{Row=9; D="1.509"; E="  -2.67%  "},
    @{Row=10; D="54.11"; E="  -2.05%  "},
    @{Row=11; D="1.005"; E="  +0.48%  "},
    @{Row=12; D="0.08907"; E="  +0.98%  "},
    @{Row=13; D="7.227"; E="  -1.04%  "},
    @{Row=14; D="23.44"; E="  +0.07%  "},
    @{Row=15; D="8.020"; E="  +5.06%  "},
    @{Row=16; D="0.00001330"; E="  -0.24%  "},
    @{Row=17; D="1.702.11"; E="  +0.08%  "},
    @{Row=18; D="100.16"; E="  -0.78%  "},
    @{Row=19; D="0.07022"; E="  -0.39%  "},
    @{Row=20; D="19.68"; E="  -0.13%  "},
    @{Row=21; D="7.022"; E="  +1.38%  "},
    @{Row=22; D=$null; E="  +0.20%  "},
    @{Row=23; D="14.51"; E="  +2.45%  "},
    @{Row=24; D="24.787.92"; E="  +0.70%  "},
    @{Row=25; D="3.228"; E="  +7.92%  "},
    @{Row=26; D="2.355"; E="  +1.23%  "},
    @{Row=27; D="22.82"; E="  +1.62%  "},
    @{Row=28; D="161.66"; E="  +1.29%  "},
    @{Row=29; D="136.45"; E="  +1.60%  "},
    @{Row=30; D="5.166"; E="  -1.07%  "},
    @{Row=31; D="7.771"; E="  -1.21%  "},
    @{Row=32; D="0.08755"; E="  +2.06%  "},
    @{Row=33; D="1.076"; E="  -3.46%  "},
    @{Row=34; D="7.190"; E="  -5.25%  "},
    @{Row=35; D="11.26"; E="  +0.40%  "},
    @{Row=36; D="1.982"; E="  +1.72%  "},
    @{Row=37; D="0.2745"; E="  -1.02%  "},
    @{Row=38; D="14.40"; E="  -3.04%  "},
    @{Row=39; D="0.09169"; E="  +1.32%  "},
    @{Row=40; D="0.02746"; E="  -1.09%  "},
    @{Row=41; D="1.465"; E="  -0.83%  "},
    @{Row=42; D="0.7694"; E="  -1.50%  "},
    @{Row=43; D="15.78"; E="  +0.67%  "},
    @{Row=44; D="0.7179"; E="  -1.89%  "},
    @{Row=45; D="2.588"; E="  +2.24%  "},
    @{Row=46; D="4.214"; E="  +0.46%  "},
    @{Row=47; D=$null; E="  +0.37%  "},
    @{Row=48; D="140.75"; E="  -0.78%  "},
    @{Row=49; D="1.314"; E="  -0.75%  "},
    @{Row=50; D="90.68"; E="  +2.67%  "},
    @{Row=51; D="0.07988"; E="  -0.59%  "}
)

foreach ($item in $updates) {
    $r = $item.Row

    if ($null -ne $item.D) {
        $cell = $ws.Range("D$r")
        # Force text storage so values like "316.97" don't get reinterpreted
        # as numbers (matches the source data, which is plain text).
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.ClearFormats()
    }

    if ($null -ne $item.E) {
        $cell = $ws.Range("E$r")
        $cell.NumberFormat = "@"
        $cell.Value = $item.E
        $cell.ClearFormats()
    }
}
